$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the date headers in B1:E1 with text labels, keeping existing
# (date) cell formatting untouched - matches typing new text over the
# old date-formatted cells in Excel.
$ws.Range("B1").Value = "nov_bef"
$ws.Range("C1").Value = "nov_aft"
$ws.Range("D1").Value = "dec_bef"
$ws.Range("E1").Value = "dec_aft"

# Update the active selection on the sheet.
$ws.Range("I8").Select()
